$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match-detail columns (F:V) between row pairs; A (index) and E (date) stay fixed per row ---
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$swapPairs = @(@(8,9), @(17,18), @(40,41))
foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $cell1 = $ws.Range($col + $r1)
        $cell2 = $ws.Range($col + $r2)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

# --- Append new match rows 44-50, copying row 43 formatting first ---
$ws.Range("A43:V43").Copy()
$ws.Range("A44:V50").PasteSpecial(-4122)

# Row 44 (Indice 43)
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "netherlands"
$ws.Range("C44").Value = "eredivisie"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45192.6875
$ws.Range("F44").Value = "FC Volendam"
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = "Heracles"
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 1.81
$ws.Range("K44").Value = "16/09/2023 20:12"
$ws.Range("L44").Value = 2.94
$ws.Range("M44").Value = "23/09/2023 16:28"
$ws.Range("N44").Value = 4.21
$ws.Range("O44").Value = "16/09/2023 20:12"
$ws.Range("P44").Value = 3.84
$ws.Range("Q44").Value = "23/09/2023 15:59"
$ws.Range("R44").Value = 4
$ws.Range("S44").Value = "16/09/2023 20:12"
$ws.Range("T44").Value = 2.36
$ws.Range("U44").Value = "23/09/2023 16:28"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/fc-volendam-heracles/M3hG12aL/"

# Row 45 (Indice 44)
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "netherlands"
$ws.Range("C45").Value = "eredivisie"
$ws.Range("D45").Value = "2023-2024"
$ws.Range("E45").Value = 45192.78125
$ws.Range("F45").Value = "Nijmegen"
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = "Utrecht"
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2.53
$ws.Range("K45").Value = "16/09/2023 20:12"
$ws.Range("L45").Value = 2.41
$ws.Range("M45").Value = "23/09/2023 18:40"
$ws.Range("N45").Value = 3.61
$ws.Range("O45").Value = "16/09/2023 20:12"
$ws.Range("P45").Value = 3.56
$ws.Range("Q45").Value = "23/09/2023 18:40"
$ws.Range("R45").Value = 2.78
$ws.Range("S45").Value = "16/09/2023 20:12"
$ws.Range("T45").Value = 3.01
$ws.Range("U45").Value = "23/09/2023 18:40"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/nijmegen-utrecht/4SsL0MER/"

# Row 46 (Indice 45)
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "netherlands"
$ws.Range("C46").Value = "eredivisie"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45192.83333333334
$ws.Range("F46").Value = "Almere City"
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = "PSV"
$ws.Range("I46").Value = 4
$ws.Range("J46").Value = 7.36
$ws.Range("K46").Value = "17/09/2023 13:43"
$ws.Range("L46").Value = 13.46
$ws.Range("M46").Value = "23/09/2023 19:59"
$ws.Range("N46").Value = 5.65
$ws.Range("O46").Value = "17/09/2023 13:43"
$ws.Range("P46").Value = 7.28
$ws.Range("Q46").Value = "23/09/2023 19:59"
$ws.Range("R46").Value = 1.38
$ws.Range("S46").Value = "17/09/2023 13:43"
$ws.Range("T46").Value = 1.21
$ws.Range("U46").Value = "23/09/2023 19:58"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-psv/tlZndtNr/"

# Row 47 (Indice 46)
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "netherlands"
$ws.Range("C47").Value = "eredivisie"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45192.875
$ws.Range("F47").Value = "Heerenveen"
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = "Excelsior"
$ws.Range("I47").Value = 3
$ws.Range("J47").Value = 1.6
$ws.Range("K47").Value = "17/09/2023 13:43"
$ws.Range("L47").Value = 1.67
$ws.Range("M47").Value = "23/09/2023 20:59"
$ws.Range("N47").Value = 4.58
$ws.Range("O47").Value = "17/09/2023 13:43"
$ws.Range("P47").Value = 4.38
$ws.Range("Q47").Value = "23/09/2023 20:59"
$ws.Range("R47").Value = 5.28
$ws.Range("S47").Value = "17/09/2023 13:43"
$ws.Range("T47").Value = 4.95
$ws.Range("U47").Value = "23/09/2023 20:59"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/heerenveen-excelsior/zBYje0xk/"

# Row 48 (Indice 47)
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "netherlands"
$ws.Range("C48").Value = "eredivisie"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45193.51041666666
$ws.Range("F48").Value = "Sparta Rotterdam"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "Vitesse"
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 2.25
$ws.Range("K48").Value = "17/09/2023 16:13"
$ws.Range("L48").Value = 2.16
$ws.Range("M48").Value = "24/09/2023 12:08"
$ws.Range("N48").Value = 3.6
$ws.Range("O48").Value = "17/09/2023 16:13"
$ws.Range("P48").Value = 3.53
$ws.Range("Q48").Value = "24/09/2023 12:14"
$ws.Range("R48").Value = 3.15
$ws.Range("S48").Value = "17/09/2023 16:13"
$ws.Range("T48").Value = 3.56
$ws.Range("U48").Value = "24/09/2023 12:14"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-vitesse/E1XffKie/"

# Row 49 (Indice 48)
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "netherlands"
$ws.Range("C49").Value = "eredivisie"
$ws.Range("D49").Value = "2023-2024"
$ws.Range("E49").Value = 45193.69791666666
$ws.Range("F49").Value = "Waalwijk"
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = "Twente"
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 4.49
$ws.Range("K49").Value = "17/09/2023 13:43"
$ws.Range("L49").Value = 6.03
$ws.Range("M49").Value = "24/09/2023 16:43"
$ws.Range("N49").Value = 4.43
$ws.Range("O49").Value = "17/09/2023 13:43"
$ws.Range("P49").Value = 4.72
$ws.Range("Q49").Value = "24/09/2023 16:44"
$ws.Range("R49").Value = 1.68
$ws.Range("S49").Value = "17/09/2023 13:43"
$ws.Range("T49").Value = 1.53
$ws.Range("U49").Value = "24/09/2023 16:38"
$ws.Range("V49").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-twente/OWEoHsa8/"

# Row 50 (Indice 49)
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "netherlands"
$ws.Range("C50").Value = "eredivisie"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45193.69791666666
$ws.Range("F50").Value = "Zwolle"
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = "AZ Alkmaar"
$ws.Range("I50").Value = 3
$ws.Range("J50").Value = 4.72
$ws.Range("K50").Value = "17/09/2023 16:13"
$ws.Range("L50").Value = 4.63
$ws.Range("M50").Value = "24/09/2023 16:39"
$ws.Range("N50").Value = 4.32
$ws.Range("O50").Value = "17/09/2023 16:13"
$ws.Range("P50").Value = 4.03
$ws.Range("Q50").Value = "24/09/2023 16:39"
$ws.Range("R50").Value = 1.67
$ws.Range("S50").Value = "17/09/2023 16:13"
$ws.Range("T50").Value = 1.77
$ws.Range("U50").Value = "24/09/2023 16:39"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-az-alkmaar/ATQ3hbM7/"

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())